# getcensus_geo_args.xlsx: switch the geography abbreviations (column B,
# "geo_abbrvs") over to the official TIGER/Line file-name abbreviations,
# per https://www2.census.gov/geo/tiger/TIGER2020/2020_TL_Shapefiles_File_Name_Definitions.pdf
# and leave the sheet scrolled/selected where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# geo_abbrvs (column B) renames -> TIGER/Line abbreviations
$ws.Range("B11").Value = "anrc"    # alaska native regional corporation:            ancsa    -> anrc
$ws.Range("B12").Value = "aiannh"  # american indian/alaska native/hawaiian home land: aiaanahhl -> aiannh
$ws.Range("B19").Value = "elsd"    # school district (elementary):                   sche     -> elsd
$ws.Range("B20").Value = "scsd"    # school district (secondary):                    schs     -> scsd
$ws.Range("B21").Value = "unsd"    # school district (unified):                      sch      -> unsd
$ws.Range("B24").Value = "sldu"    # state legislative district (upper chamber):     slupper  -> sldu
$ws.Range("B25").Value = "sldl"    # state legislative district (lower chamber):     sllower  -> sldl

# Scroll the sheet back to the top (clears the saved topLeftCell="A5") and
# leave the active selection on F11, matching the author's final view.
$excel.Goto($ws.Range("A1"), $true)
$ws.Range("F11").Select()
